$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The test data got entered in the wrong columns: column C ("Dữ Liệu Mẫu")
# held the step text and column D ("Các Bước") held the sample data.
# Fix it by swapping the C/D values for the two data rows (row 2 and 3),
# leaving the header row (row 1) untouched.

$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$ws.Range("C2").Value = $d2
$ws.Range("D2").Value = $c2

$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$ws.Range("C3").Value = $d3
$ws.Range("D3").Value = $c3

# Adjust column widths to match new (best-fit) content widths:
# column C now holds the shorter strings, column D the longer ones.
# (Input values chosen so the engine's internal pixel rounding lands as
# close as possible to the target best-fit widths of 14.02734375 / 28.45703125.)
$ws.Columns.Item(3).ColumnWidth = 13.166666666666666
$ws.Columns.Item(4).ColumnWidth = 27.666666666666668
